$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    # Assigning a numeric-looking string via .Value normally makes Excel
    # coerce it into a number (losing formatting / introducing float noise).
    # Temporarily force a text number format, set the value, then restore the
    # original cell style so no stray style attribute is left behind.
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$updates = @(
    @{ Row = 2; D = '33.836.82'; E = '  -2.16%  ' }
    @{ Row = 3; D = '1.774.58'; E = '  -0.83%  ' }
    @{ Row = 4; D = ''; E = '  -0.26%  ' }
    @{ Row = 5; D = '221.00'; E = '  -1.71%  ' }
    @{ Row = 6; D = '0.551'; E = '  -1.17%  ' }
    @{ Row = 7; D = ''; E = '  -0.08%  ' }
    @{ Row = 8; D = '31.02'; E = '  -4.60%  ' }
    @{ Row = 9; D = '0.285'; E = '  +0.53%  ' }
    @{ Row = 10; D = '0.0708'; E = '  +6.10%  ' }
    @{ Row = 11; D = ''; E = '  -1.75%  ' }
    @{ Row = 12; D = '2.029.36'; E = '  -0.88%  ' }
    @{ Row = 13; D = '1.768.84'; E = '  -1.12%  ' }
    @{ Row = 14; D = '10.52'; E = '  -4.58%  ' }
    @{ Row = 15; D = '0.624'; E = '  -1.26%  ' }
    @{ Row = 16; D = '33.830.40'; E = '  -2.27%  ' }
    @{ Row = 17; D = '4.20'; E = '  -2.11%  ' }
    @{ Row = 18; D = '67.80'; E = '  -1.56%  ' }
    @{ Row = 19; D = '244.10'; E = '  -3.85%  ' }
    @{ Row = 20; D = '0.0₃0773'; E = '  +1.08%  ' }
    @{ Row = 21; D = '0.999'; E = '  -0.18%  ' }
    @{ Row = 22; D = '10.57'; E = '  +1.66%  ' }
    @{ Row = 23; D = '4.07'; E = '  -3.73%  ' }
    @{ Row = 24; D = '2.07'; E = '  -2.60%  ' }
    @{ Row = 25; D = '157.11'; E = '  -0.98%  ' }
    @{ Row = 26; D = '16.34'; E = '  -0.32%  ' }
    @{ Row = 27; D = '6.98'; E = '  -1.54%  ' }
    @{ Row = 28; D = ''; E = '  -1.86%  ' }
    @{ Row = 29; D = ''; E = '  -0.34%  ' }
    @{ Row = 30; D = '0.0519'; E = '  +0.52%  ' }
    @{ Row = 31; D = '3.70'; E = '  -1.84%  ' }
    @{ Row = 32; D = ''; E = '  +0.53%  ' }
    @{ Row = 33; D = ''; E = '  -2.81%  ' }
    @{ Row = 34; D = ''; E = '  -2.81%  ' }
    @{ Row = 35; D = '1.393.25'; E = '  -3.48%  ' }
    @{ Row = 36; D = '0.635'; E = '  +1.51%  ' }
    @{ Row = 37; D = ''; E = '  -0.84%  ' }
    @{ Row = 38; D = ''; E = '  -1.84%  ' }
    @{ Row = 39; D = '0.929'; E = '  +3.34%  ' }
    @{ Row = 40; D = ''; E = '  -0.49%  ' }
    @{ Row = 41; D = '78.86'; E = '  -5.16%  ' }
    @{ Row = 42; D = '2.69'; E = '  -4.36%  ' }
    @{ Row = 43; D = ''; E = '  +0.44%  ' }
    @{ Row = 44; D = ''; E = '  -0.64%  ' }
    @{ Row = 45; D = '0.0489'; E = '  -2.84%  ' }
    @{ Row = 46; D = '1.03'; E = '  -1.60%  ' }
    @{ Row = 47; D = '1.926.25'; E = '  -0.95%  ' }
    @{ Row = 48; D = '104.14'; E = '  +1.11%  ' }
    @{ Row = 49; D = ''; E = '  -0.55%  ' }
    @{ Row = 50; D = '11.74'; E = '  -1.80%  ' }
)

foreach ($u in $updates) {
    if ($u.D -ne '') {
        Set-TextValue $u.Row 4 $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
